$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Лунькова Арина (row 14): grades in G and H were missing, and the other
# columns need to be raised to 5 as well.
$ws.Range("C14:H14").Value = 5

# Move the active selection (in the frozen bottom-right pane) to H14,
# matching where the user last clicked after entering the new grades.
$ws.Range("H14").Select()
